# Update Week 15 simulation totals for the "H" row on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# OFF sheet (row 2 -> "H")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 416
$wsOff.Range("C2").Value = 295
$wsOff.Range("D2").Value = 81
$wsOff.Range("E2").Value = 37

# DEF sheet (row 2 -> "H")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 423
$wsDef.Range("C2").Value = 284
$wsDef.Range("D2").Value = 101
$wsDef.Range("E2").Value = 51
$wsDef.Range("F2").Value = 7
$wsDef.Range("G2").Value = 7
